$wb = $excel.ActiveWorkbook

# Insert the new "N-wni" sheet right after "L_kier_stud" (3rd tab),
# pushing doktoranci/nauczyciele/... down by one position.
$afterSheet = $wb.Worksheets.Item(3)
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "N-wni"

# Header row
$ws.Cells.Item(1, 1).Value = "Rok"
$ws.Cells.Item(1, 2).Value = "Rodzaj"
$ws.Cells.Item(1, 3).Value = "Liczba"
$ws.Cells.Item(1, 4).Value = "Wszyscy"


$data = @(
    ,(2021, "studia stacjonarne", 385, 18333)
    ,(2021, "studia niestacjonarne", 41, 18333)
    ,(2021, "studia doktoranckie stacjonarne", 12, 392)
    ,(2021, "studia doktoranckie niestacjonarne", 3, 392)
    ,(2020, "studia stacjonarne", 419, 19993)
    ,(2020, "studia niestacjonarne", 41, 19993)
    ,(2020, "studia doktoranckie stacjonarne", 17, 566)
    ,(2020, "studia doktoranckie niestacjonarne", 4, 566)
    ,(2019, "studia stacjonarne", 434, 20692)
    ,(2019, "studia niestacjonarne", 53, 20692)
    ,(2019, "studia doktoranckie stacjonarne", 0, 725)
    ,(2019, "studia doktoranckie niestacjonarne", 0, 725)
    ,(2018, "studia stacjonarne", 434, 21509)
    ,(2018, "studia niestacjonarne", 53, 21509)
    ,(2018, "studia doktoranckie stacjonarne", 0, 954)
    ,(2018, "studia doktoranckie niestacjonarne", 0, 954)
    ,(2017, "studia stacjonarne", 497, 23410)
    ,(2017, "studia niestacjonarne", 64, 23410)
    ,(2017, "studia doktoranckie stacjonarne", 0, 919)
    ,(2017, "studia doktoranckie niestacjonarne", 0, 919)
    ,(2016, "studia stacjonarne", 504, 24398)
    ,(2016, "studia niestacjonarne", 79, 24398)
    ,(2016, "studia doktoranckie stacjonarne", 0, 952)
    ,(2016, "studia doktoranckie niestacjonarne", 0, 952)
    ,(2015, "studia stacjonarne", 483, 25423)
    ,(2015, "studia niestacjonarne", 85, 25423)
    ,(2015, "studia doktoranckie stacjonarne", 0, 953)
    ,(2015, "studia doktoranckie niestacjonarne", 0, 953)
    ,(2014, "studia stacjonarne", 523, 26974)
    ,(2014, "studia niestacjonarne", 114, 26974)
    ,(2014, "studia doktoranckie stacjonarne", 0, 917)
    ,(2014, "studia doktoranckie niestacjonarne", 0, 917)
    ,(2013, "studia stacjonarne", 539, 27959)
    ,(2013, "studia niestacjonarne", 121, 27959)
    ,(2013, "studia doktoranckie stacjonarne", 0, 881)
    ,(2013, "studia doktoranckie niestacjonarne", 0, 881)
    ,(2012, "studia stacjonarne", 585, 29802)
    ,(2012, "studia niestacjonarne", 131, 29802)
    ,(2012, "studia doktoranckie stacjonarne", 0, 876)
    ,(2012, "studia doktoranckie niestacjonarne", 0, 876)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Column B width (bestFit, matches diff's customWidth on col B)
$ws.Columns.Item(2).ColumnWidth = 34.7109375

# Selection / active cell, matching the target view state
$ws.Activate()
$ws.Range("D42").Select()

